$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Apln -> Aplnr -> ECs -------------------------------
$ws.Range("B2").Value = "Apln"
$ws.Range("C2").Value = "Aplnr"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6568066666666666
$ws.Range("H2").Value = 1.97042
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 25.63843733333334
$ws.Range("N2").Value = 76.91531200000001
$ws.Range("O2").Value = 0.9537675058687185
$ws.Range("P2").Value = 0.9537675058687185
$ws.Range("Q2").Value = 16.83949656344889
$ws.Range("R2").Value = 151.55546907104
$ws.Range("S2").Value = 0.9537675058687185
$ws.Range("T2").Value = 0.9537675058687185

# --- Row 3: ECs -> Apln -> Aplnr -> FAPs -------------------------------
$ws.Range("B3").Value = "Apln"
$ws.Range("C3").Value = "Aplnr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6568066666666666
$ws.Range("H3").Value = 1.97042
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.03172061737765654
$ws.Range("P3").Value = 0.03172061737765654
$ws.Range("Q3").Value = 0.5600518197933333
$ws.Range("R3").Value = 5.04046637814
$ws.Range("S3").Value = 0.03172061737765654
$ws.Range("T3").Value = 0.03172061737765654

# --- Row 4: ECs -> Apln -> Aplnr -> MuSCs ------------------------------
$ws.Range("B4").Value = "Apln"
$ws.Range("C4").Value = "Aplnr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6568066666666666
$ws.Range("H4").Value = 1.97042
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 0.3900970000000001
$ws.Range("N4").Value = 1.170291
$ws.Range("O4").Value = 0.01451187675362493
$ws.Range("P4").Value = 0.01451187675362493
$ws.Range("Q4").Value = 0.2562183102466667
$ws.Range("R4").Value = 2.305964792220001
$ws.Range("S4").Value = 0.01451187675362493
$ws.Range("T4").Value = 0.01451187675362493

# --- Remove the now-obsolete rows 5-7 (MuSCs sending-cluster data) ----
$ws.Range("A5:T7").Delete()
